$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look like plain numbers to Excel's parser must be
# forced to remain text (matching the original inlineStr cell type) by setting
# the number format to Text before assignment, then restoring the default style
# so no extraneous style index is left on the cell.
$forceTextCells = @("D5", "D6", "D8", "D10", "D11", "D12", "D13", "D16", "D18", "D19", "D20", "D21", "D22", "D23", "D24", "D26", "D27", "D28", "D30", "D31", "D32", "D34", "D35", "D36", "D37", "D38", "D39", "D40", "D42", "D43", "D44", "D45", "D46", "D47", "D48", "D49", "D50", "D51")
foreach ($addr in $forceTextCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# --- Coin / Link / Price / Volume cell updates ---
$ws.Range("D2").Value = '63.920.58'
$ws.Range("E2").Value = '  +1.01%  '
$ws.Range("D3").Value = '2.614.35'
$ws.Range("E3").Value = '  -0.93%  '
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("D5").Value = '597.51'
$ws.Range("E5").Value = '  -0.68%  '
$ws.Range("D6").Value = '149.73'
$ws.Range("E6").Value = '  +2.65%  '
$ws.Range("D8").Value = '0.586'
$ws.Range("E8").Value = '  +0.17%  '
$ws.Range("E9").Value = '  +1.20%  '
$ws.Range("B10").Value = 'Cardano'
$ws.Range("C10").Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range("D10").Value = '0.383'
$ws.Range("E10").Value = '  +5.57%  '
$ws.Range("B11").Value = 'Toncoin'
$ws.Range("C11").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D11").Value = '5.60'
$ws.Range("E11").Value = '  -0.12%  '
$ws.Range("D12").Value = '0.150'
$ws.Range("E12").Value = '  -0.84%  '
$ws.Range("D13").Value = '27.47'
$ws.Range("E13").Value = '  +0.98%  '
$ws.Range("D14").Value = '3.085.20'
$ws.Range("E14").Value = '  -0.81%  '
$ws.Range("D15").Value = '63.785.56'
$ws.Range("E15").Value = '  +0.96%  '
$ws.Range("D16").Value = '0.0000148'
$ws.Range("E16").Value = '  +2.64%  '
$ws.Range("D17").Value = '2.621.99'
$ws.Range("E17").Value = '  +0.08%  '
$ws.Range("D18").Value = '12.12'
$ws.Range("E18").Value = '  +6.31%  '
$ws.Range("D19").Value = '4.62'
$ws.Range("E19").Value = '  +2.99%  '
$ws.Range("D20").Value = '349.16'
$ws.Range("E20").Value = '  +2.40%  '
$ws.Range("D21").Value = '6.87'
$ws.Range("E21").Value = '  -0.24%  '
$ws.Range("D22").Value = '0.999'
$ws.Range("E22").Value = '  -0.09%  '
$ws.Range("D23").Value = '5.71'
$ws.Range("E23").Value = '  +2.42%  '
$ws.Range("D24").Value = '66.13'
$ws.Range("E24").Value = '  -0.45%  '
$ws.Range("E25").Value = '  +14.56%  '
$ws.Range("D26").Value = '9.19'
$ws.Range("E26").Value = '  +5.58%  '
$ws.Range("D27").Value = '1.68'
$ws.Range("E27").Value = '  +2.09%  '
$ws.Range("D28").Value = '8.10'
$ws.Range("E28").Value = '  +3.61%  '
$ws.Range("E29").Value = '  +0.30%  '
$ws.Range("B30").Value = 'Bittensor'
$ws.Range("C30").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D30").Value = '542.20'
$ws.Range("E30").Value = '  -1.63%  '
$ws.Range("B31").Value = 'Binance-PegBSC-USD'
$ws.Range("C31").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D31").Value = '1.00'
$ws.Range("E31").Value = '  +0.06%  '
$ws.Range("D32").Value = '2.04'
$ws.Range("E32").Value = '  +0.35%  '
$ws.Range("D33").Value = '0.0₃0847'
$ws.Range("E33").Value = '  +5.40%  '
$ws.Range("D34").Value = '1.74'
$ws.Range("E34").Value = '  +0.09%  '
$ws.Range("D35").Value = '5.25'
$ws.Range("E35").Value = '  +0.59%  '
$ws.Range("D36").Value = '168.38'
$ws.Range("D37").Value = '0.407'
$ws.Range("E37").Value = '  +0.45%  '
$ws.Range("D38").Value = '0.999'
$ws.Range("E38").Value = '  -0.12%  '
$ws.Range("D39").Value = '1.96'
$ws.Range("E39").Value = '  +4.85%  '
$ws.Range("D40").Value = '19.39'
$ws.Range("E40").Value = '  +2.17%  '
$ws.Range("E41").Value = '  +0.07%  '
$ws.Range("D42").Value = '167.43'
$ws.Range("E42").Value = '  -0.28%  '
$ws.Range("D43").Value = '39.86'
$ws.Range("E43").Value = '  +0.59%  '
$ws.Range("D44").Value = '3.91'
$ws.Range("E44").Value = '  +4.65%  '
$ws.Range("D45").Value = '0.0585'
$ws.Range("E45").Value = '  +2.49%  '
$ws.Range("D46").Value = '21.34'
$ws.Range("E46").Value = '  -4.91%  '
$ws.Range("D47").Value = '0.625'
$ws.Range("E47").Value = '  +0.33%  '
$ws.Range("D48").Value = '1.98'
$ws.Range("E48").Value = '  +12.96%  '
$ws.Range("D49").Value = '0.0244'
$ws.Range("E49").Value = '  +0.52%  '
$ws.Range("D50").Value = '0.0963'
$ws.Range("E50").Value = '  +0.46%  '
$ws.Range("D51").Value = '19.08'
$ws.Range("E51").Value = '  +2.23%  '

# Restore default (unstyled) formatting on cells where we forced Text format,
# so the resulting style matches the workbook's original (unstyled) cells.
foreach ($addr in $forceTextCells) {
    $ws.Range($addr).Style = "Normal"
}
